# "Add files via upload" — OSAV_T3_GTFS_01022020.pptx
#
# Content-level change captured by the diff: on the "Dev Environments"
# slide (slide 6), the bullet that read "Python 3.7" was retyped as
# "Python (3.7)".
#
# (The diff's other hunks — embeddedFontLst bookkeeping in
# presentation.xml and the scattered rPr/endParaRPr dirty="0" markers —
# are PowerPoint-desktop resave artifacts with no corresponding
# property/method on the Slide/Shape/TextRange object model, so they
# aren't reproducible through COM automation; only the genuine text
# edit below is applied.)

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(6)
$shape = $slide.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

$textRange.Replace("Python 3.7", "Python (3.7)") | Out-Null
